$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2134292565947242
$ws.Cells.Item(2, 3).Value = 0.5251798561151079
$ws.Cells.Item(2, 10).Value = 0.02398081534772182
$ws.Cells.Item(2, 16).Value = 0.1654676258992806
$ws.Cells.Item(2, 19).Value = 0.07194244604316546
$ws.Cells.Item(3, 2).Value = 0.008547008547008548
$ws.Cells.Item(3, 3).Value = 0.0170940170940171
$ws.Cells.Item(3, 10).Value = 0.04273504273504274
$ws.Cells.Item(3, 15).Value = 0.008547008547008548
$ws.Cells.Item(3, 16).Value = 0.7350427350427351
$ws.Cells.Item(3, 19).Value = 0.188034188034188
$ws.Cells.Item(4, 10).Value = 0.1
$ws.Cells.Item(4, 16).Value = 0.58
$ws.Cells.Item(4, 19).Value = 0.32
$ws.Cells.Item(6, 2).Value = 0.04
$ws.Cells.Item(6, 4).Value = 0.02222222222222222
$ws.Cells.Item(6, 6).Value = 0.05333333333333334
$ws.Cells.Item(6, 10).Value = 0.2888888888888889
$ws.Cells.Item(6, 15).Value = 0.02666666666666667
$ws.Cells.Item(6, 17).Value = 0.1688888888888889
$ws.Cells.Item(6, 18).Value = 0.06222222222222222
$ws.Cells.Item(6, 19).Value = 0.3377777777777778
$ws.Cells.Item(7, 2).Value = 0.1420118343195266
$ws.Cells.Item(7, 4).Value = 0.005917159763313609
$ws.Cells.Item(7, 6).Value = 0.0650887573964497
$ws.Cells.Item(7, 10).Value = 0.1479289940828402
$ws.Cells.Item(7, 15).Value = 0.01775147928994083
$ws.Cells.Item(7, 17).Value = 0.1775147928994083
$ws.Cells.Item(7, 18).Value = 0.0650887573964497
$ws.Cells.Item(7, 19).Value = 0.378698224852071
$ws.Cells.Item(8, 2).Value = 0.1294642857142857
$ws.Cells.Item(8, 4).Value = 0.01785714285714286
$ws.Cells.Item(8, 5).Value = 0.002232142857142857
$ws.Cells.Item(8, 6).Value = 0.05133928571428571
$ws.Cells.Item(8, 10).Value = 0.09821428571428571
$ws.Cells.Item(8, 15).Value = 0.02008928571428572
$ws.Cells.Item(8, 17).Value = 0.2053571428571428
$ws.Cells.Item(8, 18).Value = 0.08705357142857142
$ws.Cells.Item(8, 19).Value = 0.3883928571428572
$ws.Cells.Item(9, 2).Value = 0.142156862745098
$ws.Cells.Item(9, 4).Value = 0.02450980392156863
$ws.Cells.Item(9, 6).Value = 0.09313725490196079
$ws.Cells.Item(9, 10).Value = 0.1127450980392157
$ws.Cells.Item(9, 15).Value = 0.02450980392156863
$ws.Cells.Item(9, 17).Value = 0.1813725490196078
$ws.Cells.Item(9, 18).Value = 0.06862745098039216
$ws.Cells.Item(9, 19).Value = 0.3529411764705883
$ws.Cells.Item(10, 2).Value = 0.1384505021520804
$ws.Cells.Item(10, 4).Value = 0.02295552367288379
$ws.Cells.Item(10, 5).Value = 0.0007173601147776184
$ws.Cells.Item(10, 6).Value = 0.06169296987087518
$ws.Cells.Item(10, 10).Value = 0.1327116212338594
$ws.Cells.Item(10, 15).Value = 0.02223816355810617
$ws.Cells.Item(10, 17).Value = 0.1951219512195122
$ws.Cells.Item(10, 18).Value = 0.07101865136298421
$ws.Cells.Item(10, 19).Value = 0.3550932568149211
$ws.Cells.Item(11, 7).Value = 0.1615120274914089
$ws.Cells.Item(11, 10).Value = 0.1237113402061856
$ws.Cells.Item(11, 11).Value = 0.2646048109965636
$ws.Cells.Item(11, 12).Value = 0.4398625429553265
$ws.Cells.Item(11, 19).Value = 0.01030927835051546
$ws.Cells.Item(12, 7).Value = 0.6438356164383562
$ws.Cells.Item(12, 10).Value = 0.1986301369863014
$ws.Cells.Item(12, 11).Value = 0.00684931506849315
$ws.Cells.Item(12, 12).Value = 0.1095890410958904
$ws.Cells.Item(12, 19).Value = 0.0410958904109589
$ws.Cells.Item(13, 7).Value = 0.5918367346938775
$ws.Cells.Item(13, 10).Value = 0.3469387755102041
$ws.Cells.Item(13, 19).Value = 0.06122448979591837
$ws.Cells.Item(14, 7).Value = 0.8333333333333334
$ws.Cells.Item(14, 10).Value = 0.1666666666666667
$ws.Cells.Item(15, 6).Value = 0.004273504273504274
$ws.Cells.Item(15, 8).Value = 0.1452991452991453
$ws.Cells.Item(15, 9).Value = 0.07692307692307693
$ws.Cells.Item(15, 10).Value = 0.3803418803418803
$ws.Cells.Item(15, 11).Value = 0.05555555555555555
$ws.Cells.Item(15, 13).Value = 0.004273504273504274
$ws.Cells.Item(15, 15).Value = 0.03418803418803419
$ws.Cells.Item(15, 19).Value = 0.2991452991452991
$ws.Cells.Item(16, 8).Value = 0.2068965517241379
$ws.Cells.Item(16, 9).Value = 0.09961685823754789
$ws.Cells.Item(16, 10).Value = 0.4061302681992337
$ws.Cells.Item(16, 11).Value = 0.08812260536398467
$ws.Cells.Item(16, 13).Value = 0.02298850574712644
$ws.Cells.Item(16, 14).Value = 0.01532567049808429
$ws.Cells.Item(16, 15).Value = 0.05363984674329502
$ws.Cells.Item(16, 19).Value = 0.10727969348659
$ws.Cells.Item(17, 6).Value = 0.02330508474576271
$ws.Cells.Item(17, 8).Value = 0.1504237288135593
$ws.Cells.Item(17, 9).Value = 0.08898305084745763
$ws.Cells.Item(17, 10).Value = 0.4470338983050847
$ws.Cells.Item(17, 11).Value = 0.09110169491525423
$ws.Cells.Item(17, 13).Value = 0.0211864406779661
$ws.Cells.Item(17, 15).Value = 0.06567796610169492
$ws.Cells.Item(17, 19).Value = 0.1122881355932203
$ws.Cells.Item(18, 6).Value = 0.01111111111111111
$ws.Cells.Item(18, 8).Value = 0.1388888888888889
$ws.Cells.Item(18, 9).Value = 0.1055555555555556
$ws.Cells.Item(18, 10).Value = 0.4444444444444444
$ws.Cells.Item(18, 11).Value = 0.08888888888888889
$ws.Cells.Item(18, 13).Value = 0.01111111111111111
$ws.Cells.Item(18, 14).Value = 0.01111111111111111
$ws.Cells.Item(18, 15).Value = 0.06666666666666667
$ws.Cells.Item(18, 19).Value = 0.1222222222222222
$ws.Cells.Item(19, 6).Value = 0.01825396825396826
$ws.Cells.Item(19, 8).Value = 0.2079365079365079
$ws.Cells.Item(19, 9).Value = 0.08492063492063492
$ws.Cells.Item(19, 10).Value = 0.3833333333333334
$ws.Cells.Item(19, 11).Value = 0.09841269841269841
$ws.Cells.Item(19, 13).Value = 0.02698412698412699
$ws.Cells.Item(19, 14).Value = 0.001587301587301587
$ws.Cells.Item(19, 15).Value = 0.08015873015873017
$ws.Cells.Item(19, 19).Value = 0.09841269841269841

Write-Output "Applied 113 cell updates"
